$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting existing rows 36..263 down to 37..264
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly data point
$ws.Cells.Item(36, 1).Value = 3
$ws.Cells.Item(36, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 44602
$ws.Cells.Item(36, 5).Value = 5
$ws.Cells.Item(36, 6).Value = 100112039
$ws.Cells.Item(36, 7).Value = "Ciboulette"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 130
$ws.Cells.Item(36, 11).Value = 1500
$ws.Cells.Item(36, 12).Value = 1500
$ws.Cells.Item(36, 13).Value = 1500
$ws.Cells.Item(36, 14).Value = "`$/docena de atados"
$ws.Cells.Item(36, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(36, 16).Value = 500
$ws.Cells.Item(36, 17).Value = 3
$ws.Cells.Item(36, 18).Value = "Hortaliza"
